# Apply the "Finished the equipment section" update to the Items sheet.
# - Adds 36 new equipment/spell/artifact rows (97-132) with craft info.
# - Adds craft XP bonus columns (U-X) to existing row 93.
# - Widens columns F and K slightly to fit the new, longer content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

$ws.Cells.Item(93,21).Value = 1
$ws.Cells.Item(93,22).Value = 9
$ws.Cells.Item(93,23).Value = 22
$ws.Cells.Item(93,24).Value = 'spell'
$ws.Cells.Item(97,1).Value = 'Blood Lust'
$ws.Cells.Item(97,3).Value = 1
$ws.Cells.Item(97,4).Value = 'Sharp Knife'
$ws.Cells.Item(97,5).Value = 'weapon'
$ws.Cells.Item(97,6).Value = 'A very sharp knife, careful you don''t cut your self with that.'
$ws.Cells.Item(97,8).Value = 5
$ws.Cells.Item(97,11).Value = 25
$ws.Cells.Item(97,15).Value = 0.05
$ws.Cells.Item(97,16).Value = 0.05
$ws.Cells.Item(97,17).Value = 0.05
$ws.Cells.Item(97,18).Value = 0.05
$ws.Cells.Item(97,19).Value = 0.05
$ws.Cells.Item(97,21).Value = 1
$ws.Cells.Item(97,22).Value = 2
$ws.Cells.Item(97,23).Value = 5
$ws.Cells.Item(97,24).Value = 'weapon'
$ws.Cells.Item(98,1).Value = 'Prayer Of Faith'
$ws.Cells.Item(98,3).Value = 1
$ws.Cells.Item(98,4).Value = 'Cloth Shoes'
$ws.Cells.Item(98,5).Value = 'feet'
$ws.Cells.Item(98,6).Value = 'Simple shoes to keep your feet warm.'
$ws.Cells.Item(98,7).Value = 'feet'
$ws.Cells.Item(98,10).Value = 3
$ws.Cells.Item(98,11).Value = 40
$ws.Cells.Item(98,21).Value = 1
$ws.Cells.Item(98,22).Value = 2
$ws.Cells.Item(98,23).Value = 5
$ws.Cells.Item(98,24).Value = 'armour'
$ws.Cells.Item(99,1).Value = 'Blood Lust'
$ws.Cells.Item(99,2).Value = 'Artifact Hunter'
$ws.Cells.Item(99,3).Value = 1
$ws.Cells.Item(99,4).Value = 'Fingerless ripped gloves'
$ws.Cells.Item(99,5).Value = 'gloves'
$ws.Cells.Item(99,6).Value = 'These fingerless ripped gloves have smelt better...At one time.'
$ws.Cells.Item(99,7).Value = 'hands'
$ws.Cells.Item(99,10).Value = 1
$ws.Cells.Item(99,11).Value = 10
$ws.Cells.Item(99,21).Value = 1
$ws.Cells.Item(99,22).Value = 0
$ws.Cells.Item(99,23).Value = 5
$ws.Cells.Item(99,24).Value = 'armour'
$ws.Cells.Item(100,1).Value = 'Prayer Of Faith'
$ws.Cells.Item(100,2).Value = 'Enchantresses Heart'
$ws.Cells.Item(100,3).Value = 1
$ws.Cells.Item(100,4).Value = 'Slayers Gloves'
$ws.Cells.Item(100,5).Value = 'gloves'
$ws.Cells.Item(100,6).Value = 'These gloves are worn by the slayers of hell who hunt the enemies with no names.'
$ws.Cells.Item(100,7).Value = 'gloves'
$ws.Cells.Item(100,10).Value = 17
$ws.Cells.Item(100,11).Value = 160
$ws.Cells.Item(100,15).Value = 0.15
$ws.Cells.Item(100,16).Value = 0.13
$ws.Cells.Item(100,17).Value = 0.12
$ws.Cells.Item(100,18).Value = 0.1
$ws.Cells.Item(100,19).Value = 0.11
$ws.Cells.Item(100,21).Value = 1
$ws.Cells.Item(100,22).Value = 13
$ws.Cells.Item(100,23).Value = 30
$ws.Cells.Item(100,24).Value = 'armour'
$ws.Cells.Item(101,1).Value = 'Squires Hopes'
$ws.Cells.Item(101,3).Value = 1
$ws.Cells.Item(101,4).Value = 'Long Sword'
$ws.Cells.Item(101,5).Value = 'weapon'
$ws.Cells.Item(101,6).Value = 'Simple but effective for a sword of it''s type.'
$ws.Cells.Item(101,8).Value = 8
$ws.Cells.Item(101,11).Value = 50
$ws.Cells.Item(101,15).Value = 0.08
$ws.Cells.Item(101,16).Value = 0.08
$ws.Cells.Item(101,17).Value = 0.08
$ws.Cells.Item(101,18).Value = 0.08
$ws.Cells.Item(101,19).Value = 0.08
$ws.Cells.Item(101,21).Value = 1
$ws.Cells.Item(101,22).Value = 4
$ws.Cells.Item(101,23).Value = 8
$ws.Cells.Item(101,24).Value = 'weapon'
$ws.Cells.Item(102,2).Value = 'Smiths Blood'
$ws.Cells.Item(102,3).Value = 1
$ws.Cells.Item(102,4).Value = 'Demon Scale Leggings'
$ws.Cells.Item(102,5).Value = 'leggings'
$ws.Cells.Item(102,6).Value = 'Created from the scales of demons, these leggings are cursed by a thousand sins.'
$ws.Cells.Item(102,7).Value = 'leggings'
$ws.Cells.Item(102,10).Value = 18
$ws.Cells.Item(102,11).Value = 130
$ws.Cells.Item(102,15).Value = 0.1
$ws.Cells.Item(102,16).Value = 0.13
$ws.Cells.Item(102,17).Value = 0.1
$ws.Cells.Item(102,18).Value = 0.1
$ws.Cells.Item(102,19).Value = 0.13
$ws.Cells.Item(102,21).Value = 1
$ws.Cells.Item(102,22).Value = 9
$ws.Cells.Item(102,23).Value = 30
$ws.Cells.Item(102,24).Value = 'armour'
$ws.Cells.Item(103,1).Value = 'Blood Lust'
$ws.Cells.Item(103,3).Value = 1
$ws.Cells.Item(103,4).Value = 'Wooden Round Shield'
$ws.Cells.Item(103,5).Value = 'shield'
$ws.Cells.Item(103,6).Value = 'A simple, well carved and well enforced wooden shield.'
$ws.Cells.Item(103,10).Value = 3
$ws.Cells.Item(103,11).Value = 25
$ws.Cells.Item(103,21).Value = 1
$ws.Cells.Item(103,22).Value = 2
$ws.Cells.Item(103,23).Value = 6
$ws.Cells.Item(103,24).Value = 'armour'
$ws.Cells.Item(104,1).Value = 'Archers Bane'
$ws.Cells.Item(104,3).Value = 1
$ws.Cells.Item(104,4).Value = 'Bloody Ice Sleeves'
$ws.Cells.Item(104,5).Value = 'sleeves'
$ws.Cells.Item(104,6).Value = 'Sleeves made of the finest crystal in the shape of ice crawling up your arms, dripping in blood.'
$ws.Cells.Item(104,7).Value = 'sleeves'
$ws.Cells.Item(104,10).Value = 16
$ws.Cells.Item(104,11).Value = 700
$ws.Cells.Item(104,15).Value = 0.16
$ws.Cells.Item(104,16).Value = 0.18
$ws.Cells.Item(104,17).Value = 0.16
$ws.Cells.Item(104,21).Value = 1
$ws.Cells.Item(104,22).Value = 29
$ws.Cells.Item(104,23).Value = 60
$ws.Cells.Item(104,24).Value = 'armour'
$ws.Cells.Item(105,1).Value = 'Goblin War Cry'
$ws.Cells.Item(105,3).Value = 1
$ws.Cells.Item(105,4).Value = 'Rusty bloody broken dagger'
$ws.Cells.Item(105,5).Value = 'weapon'
$ws.Cells.Item(105,8).Value = 3
$ws.Cells.Item(105,11).Value = 10
$ws.Cells.Item(105,21).Value = 1
$ws.Cells.Item(105,22).Value = 0
$ws.Cells.Item(105,23).Value = 5
$ws.Cells.Item(105,24).Value = 'weapon'
$ws.Cells.Item(106,4).Value = 'Litch Scythe'
$ws.Cells.Item(106,5).Value = 'weapon'
$ws.Cells.Item(106,6).Value = 'The scythe of a list carries the powers of the dead, of death to be more specific. Made of the darkest magics and held together by the sacrifices of the damned, this scythe will cut through both light and shadow.'
$ws.Cells.Item(106,8).Value = 50
$ws.Cells.Item(106,11).Value = 750
$ws.Cells.Item(106,15).Value = 0.19
$ws.Cells.Item(106,16).Value = 0.19
$ws.Cells.Item(106,17).Value = 0.19
$ws.Cells.Item(106,18).Value = 0.19
$ws.Cells.Item(106,19).Value = 0.19
$ws.Cells.Item(106,21).Value = 1
$ws.Cells.Item(106,22).Value = 13
$ws.Cells.Item(106,23).Value = 25
$ws.Cells.Item(106,24).Value = 'weapon'
$ws.Cells.Item(107,2).Value = 'Golden Touch'
$ws.Cells.Item(107,3).Value = 1
$ws.Cells.Item(107,4).Value = 'Leather Pants'
$ws.Cells.Item(107,5).Value = 'leggings'
$ws.Cells.Item(107,6).Value = 'Simple leather pants.'
$ws.Cells.Item(107,7).Value = 'leggings'
$ws.Cells.Item(107,10).Value = 3
$ws.Cells.Item(107,11).Value = 50
$ws.Cells.Item(107,21).Value = 1
$ws.Cells.Item(107,22).Value = 2
$ws.Cells.Item(107,23).Value = 5
$ws.Cells.Item(107,24).Value = 'armour'
$ws.Cells.Item(108,4).Value = 'Dagger Of Mystics'
$ws.Cells.Item(108,5).Value = 'weapon'
$ws.Cells.Item(108,6).Value = 'Theres a dagger that acts like a key. A key to a mystical plane known as the astral plane. How ever, much like this dagger, the mystics also made this dagger.'
$ws.Cells.Item(108,8).Value = 70
$ws.Cells.Item(108,11).Value = 1800
$ws.Cells.Item(108,15).Value = 0.2
$ws.Cells.Item(108,16).Value = 0.21
$ws.Cells.Item(108,17).Value = 0.22
$ws.Cells.Item(108,18).Value = 0.23
$ws.Cells.Item(108,19).Value = 0.24
$ws.Cells.Item(108,21).Value = 1
$ws.Cells.Item(108,22).Value = 15
$ws.Cells.Item(108,23).Value = 30
$ws.Cells.Item(108,24).Value = 'weapon'
$ws.Cells.Item(109,4).Value = 'Witches Hooked Broom'
$ws.Cells.Item(109,5).Value = 'weapon'
$ws.Cells.Item(109,6).Value = 'This hooked broom can let you fly high into the sky, and rip out their cuts with the razor sharp hook on the end of the handle. After wards, you can use the broom to clean up said guts.'
$ws.Cells.Item(109,8).Value = 120
$ws.Cells.Item(109,11).Value = 2590
$ws.Cells.Item(109,15).Value = 0.23
$ws.Cells.Item(109,16).Value = 0.25
$ws.Cells.Item(109,17).Value = 0.24
$ws.Cells.Item(109,18).Value = 0.25
$ws.Cells.Item(109,19).Value = 0.25
$ws.Cells.Item(109,21).Value = 1
$ws.Cells.Item(109,22).Value = 18
$ws.Cells.Item(109,23).Value = 35
$ws.Cells.Item(109,24).Value = 'weapon'
$ws.Cells.Item(110,4).Value = 'Dwarven Forged Breast'
$ws.Cells.Item(110,5).Value = 'body'
$ws.Cells.Item(110,6).Value = 'Forged by dwarves in their mountain homes, deep in the pits of the earth when you people like to go and hunt for shiny objects. More like get your selves killed while looking for those "shiny pieces".'
$ws.Cells.Item(110,7).Value = 'body'
$ws.Cells.Item(110,10).Value = 150
$ws.Cells.Item(110,11).Value = 5000
$ws.Cells.Item(110,15).Value = 0.29
$ws.Cells.Item(110,16).Value = 0.29
$ws.Cells.Item(110,17).Value = 0.29
$ws.Cells.Item(110,18).Value = 0.29
$ws.Cells.Item(110,19).Value = 0.29
$ws.Cells.Item(110,21).Value = 1
$ws.Cells.Item(110,22).Value = 40
$ws.Cells.Item(110,23).Value = 70
$ws.Cells.Item(110,24).Value = 'armour'
$ws.Cells.Item(111,4).Value = 'Witches Flesh Shield'
$ws.Cells.Item(111,5).Value = 'shield'
$ws.Cells.Item(111,6).Value = 'Made from the the finest and freshest witch flesh I could find. Skinned her alive I did.'
$ws.Cells.Item(111,10).Value = 50
$ws.Cells.Item(111,11).Value = 3000
$ws.Cells.Item(111,15).Value = 0.2
$ws.Cells.Item(111,16).Value = 0.2
$ws.Cells.Item(111,17).Value = 0.2
$ws.Cells.Item(111,18).Value = 0.2
$ws.Cells.Item(111,19).Value = 0.2
$ws.Cells.Item(111,21).Value = 1
$ws.Cells.Item(111,22).Value = 30
$ws.Cells.Item(111,23).Value = 65
$ws.Cells.Item(111,24).Value = 'armour'
$ws.Cells.Item(112,4).Value = 'Elven Chain'
$ws.Cells.Item(112,5).Value = 'leggings'
$ws.Cells.Item(112,6).Value = 'Th elves are known for creating some very beautiful art, weapons and most of all armour.'
$ws.Cells.Item(112,7).Value = 'leggings'
$ws.Cells.Item(112,10).Value = 25
$ws.Cells.Item(112,11).Value = 1200
$ws.Cells.Item(112,15).Value = 0.17
$ws.Cells.Item(112,16).Value = 0.17
$ws.Cells.Item(112,17).Value = 0.17
$ws.Cells.Item(112,18).Value = 0.18
$ws.Cells.Item(112,19).Value = 0.18
$ws.Cells.Item(112,21).Value = 1
$ws.Cells.Item(112,22).Value = 29
$ws.Cells.Item(112,23).Value = 65
$ws.Cells.Item(112,24).Value = 'armour'
$ws.Cells.Item(113,4).Value = 'Cursed Iron Cap'
$ws.Cells.Item(113,5).Value = 'helmet'
$ws.Cells.Item(113,6).Value = 'Cursed by some priest in some old village, this cap will fill your head with the voices of the damned.'
$ws.Cells.Item(113,10).Value = 22
$ws.Cells.Item(113,11).Value = 1657
$ws.Cells.Item(113,15).Value = 0.18
$ws.Cells.Item(113,16).Value = 0.22
$ws.Cells.Item(113,17).Value = 0.18
$ws.Cells.Item(113,18).Value = 0.18
$ws.Cells.Item(113,19).Value = 0.19
$ws.Cells.Item(113,21).Value = 1
$ws.Cells.Item(113,22).Value = 32
$ws.Cells.Item(113,23).Value = 67
$ws.Cells.Item(113,24).Value = 'armour'
$ws.Cells.Item(114,4).Value = 'Devils Hand Shake'
$ws.Cells.Item(114,5).Value = 'gloves'
$ws.Cells.Item(114,6).Value = 'Simple gloves, they are. But imbued with the hand shake of the devil him self. They say theres good luck with in gloves.'
$ws.Cells.Item(114,7).Value = 'gloves'
$ws.Cells.Item(114,10).Value = 18
$ws.Cells.Item(114,11).Value = 1600
$ws.Cells.Item(114,15).Value = 0.2
$ws.Cells.Item(114,16).Value = 0.22
$ws.Cells.Item(114,17).Value = 0.2
$ws.Cells.Item(114,18).Value = 0.22
$ws.Cells.Item(114,19).Value = 0.2
$ws.Cells.Item(114,21).Value = 1
$ws.Cells.Item(114,22).Value = 34
$ws.Cells.Item(114,23).Value = 67
$ws.Cells.Item(114,24).Value = 'armour'
$ws.Cells.Item(115,4).Value = 'Angelic Plate Boots'
$ws.Cells.Item(115,5).Value = 'feet'
$ws.Cells.Item(115,6).Value = 'Blessed by the angels them selves, these boots will carry you high into the heavens.'
$ws.Cells.Item(115,7).Value = 'feet'
$ws.Cells.Item(115,10).Value = 18
$ws.Cells.Item(115,11).Value = 1300
$ws.Cells.Item(115,15).Value = 0.22
$ws.Cells.Item(115,16).Value = 0.23
$ws.Cells.Item(115,17).Value = 0.24
$ws.Cells.Item(115,18).Value = 0.24
$ws.Cells.Item(115,19).Value = 0.2
$ws.Cells.Item(115,21).Value = 1
$ws.Cells.Item(115,22).Value = 28
$ws.Cells.Item(115,23).Value = 65
$ws.Cells.Item(115,24).Value = 'armour'
$ws.Cells.Item(116,4).Value = 'Goblin Leather Sleeves'
$ws.Cells.Item(116,5).Value = 'sleeves'
$ws.Cells.Item(116,6).Value = 'Made from the flesh of goblins, stitched with their hair and held together with hope.'
$ws.Cells.Item(116,10).Value = 18
$ws.Cells.Item(116,11).Value = 1800
$ws.Cells.Item(116,16).Value = 0.22
$ws.Cells.Item(116,17).Value = 0.22
$ws.Cells.Item(116,18).Value = 0.21
$ws.Cells.Item(116,21).Value = 1
$ws.Cells.Item(116,22).Value = 31
$ws.Cells.Item(116,23).Value = 68
$ws.Cells.Item(116,24).Value = 'armour'
$ws.Cells.Item(117,4).Value = 'Shadows Dance'
$ws.Cells.Item(117,5).Value = 'spell-damage'
$ws.Cells.Item(117,6).Value = 'Make the shadows dance for you, make them hunt for you. Make them kill for you.'
$ws.Cells.Item(117,8).Value = 160
$ws.Cells.Item(117,11).Value = 3200
$ws.Cells.Item(117,21).Value = 1
$ws.Cells.Item(117,22).Value = 12
$ws.Cells.Item(117,23).Value = 25
$ws.Cells.Item(117,24).Value = 'spell'
$ws.Cells.Item(118,4).Value = 'Unholy Vow'
$ws.Cells.Item(118,5).Value = 'spell-healing'
$ws.Cells.Item(118,6).Value = 'Take the vow and become a sinful, prideful person - as long as you take the vow child, your wounds shall heal over.'
$ws.Cells.Item(118,9).Value = 120
$ws.Cells.Item(118,11).Value = 2400
$ws.Cells.Item(118,21).Value = 1
$ws.Cells.Item(118,22).Value = 14
$ws.Cells.Item(118,23).Value = 30
$ws.Cells.Item(118,24).Value = 'spell'
$ws.Cells.Item(119,4).Value = 'Astral Ring'
$ws.Cells.Item(119,5).Value = 'ring'
$ws.Cells.Item(119,6).Value = 'Looking to head to the astral plane? This isn''t the key the mystics made, but it will focus your mind and soul.'
$ws.Cells.Item(119,8).Value = 150
$ws.Cells.Item(119,11).Value = 2500
$ws.Cells.Item(119,15).Value = 0.21
$ws.Cells.Item(119,16).Value = 0.22
$ws.Cells.Item(119,17).Value = 0.22
$ws.Cells.Item(119,18).Value = 0.25
$ws.Cells.Item(119,19).Value = 0.25
$ws.Cells.Item(119,21).Value = 1
$ws.Cells.Item(119,22).Value = 23
$ws.Cells.Item(119,23).Value = 55
$ws.Cells.Item(119,24).Value = 'ring'
$ws.Cells.Item(120,1).Value = 'Prayer Of Faith'
$ws.Cells.Item(120,3).Value = 1
$ws.Cells.Item(120,4).Value = 'Wizards Hat'
$ws.Cells.Item(120,5).Value = 'helmet'
$ws.Cells.Item(120,6).Value = 'Gives and unfair advantage to those of the magical kind.'
$ws.Cells.Item(120,7).Value = 'helmet'
$ws.Cells.Item(120,10).Value = 16
$ws.Cells.Item(120,11).Value = 270
$ws.Cells.Item(120,15).Value = 0.12
$ws.Cells.Item(120,16).Value = 0.12
$ws.Cells.Item(120,17).Value = 0.12
$ws.Cells.Item(120,18).Value = 0.16
$ws.Cells.Item(120,19).Value = 0.16
$ws.Cells.Item(120,21).Value = 1
$ws.Cells.Item(120,22).Value = 20
$ws.Cells.Item(120,23).Value = 50
$ws.Cells.Item(120,24).Value = 'armour'
$ws.Cells.Item(121,4).Value = 'Dragon Mage Scale'
$ws.Cells.Item(121,5).Value = 'body'
$ws.Cells.Item(121,6).Value = 'Created from the scales of dragons and enchanted by the mages of old, this armour is hard to come by.'
$ws.Cells.Item(121,7).Value = 'body'
$ws.Cells.Item(121,10).Value = 200
$ws.Cells.Item(121,11).Value = 10000
$ws.Cells.Item(121,15).Value = 0.33
$ws.Cells.Item(121,16).Value = 0.33
$ws.Cells.Item(121,17).Value = 0.33
$ws.Cells.Item(121,18).Value = 0.33
$ws.Cells.Item(121,19).Value = 0.33
$ws.Cells.Item(121,21).Value = 1
$ws.Cells.Item(121,22).Value = 48
$ws.Cells.Item(121,23).Value = 70
$ws.Cells.Item(121,24).Value = 'armour'
$ws.Cells.Item(122,4).Value = 'Ice Fiend Shield'
$ws.Cells.Item(122,5).Value = 'shield'
$ws.Cells.Item(122,6).Value = 'I hunted these beasts, day in and day out. Crafted weapons and armour from their corpses. Alas this is my proudest possession.'
$ws.Cells.Item(122,10).Value = 75
$ws.Cells.Item(122,11).Value = 5400
$ws.Cells.Item(122,15).Value = 0.23
$ws.Cells.Item(122,16).Value = 0.24
$ws.Cells.Item(122,17).Value = 0.23
$ws.Cells.Item(122,18).Value = 0.24
$ws.Cells.Item(122,19).Value = 0.23
$ws.Cells.Item(122,21).Value = 1
$ws.Cells.Item(122,22).Value = 34
$ws.Cells.Item(122,23).Value = 75
$ws.Cells.Item(122,24).Value = 'armour'
$ws.Cells.Item(123,4).Value = 'Blood Covered Stone Leggings'
$ws.Cells.Item(123,5).Value = 'leggings'
$ws.Cells.Item(123,6).Value = 'Enchanted to look like stone, but as light as a feather. These leggings are covered in the blood of another how ever.'
$ws.Cells.Item(123,7).Value = 'leggings'
$ws.Cells.Item(123,10).Value = 28
$ws.Cells.Item(123,11).Value = 2600
$ws.Cells.Item(123,15).Value = 0.19
$ws.Cells.Item(123,16).Value = 0.22
$ws.Cells.Item(123,17).Value = 0.19
$ws.Cells.Item(123,18).Value = 0.23
$ws.Cells.Item(123,19).Value = 0.22
$ws.Cells.Item(123,21).Value = 1
$ws.Cells.Item(123,22).Value = 36
$ws.Cells.Item(123,23).Value = 78
$ws.Cells.Item(123,24).Value = 'armour'
$ws.Cells.Item(124,4).Value = 'Ageless Leather Boots'
$ws.Cells.Item(124,5).Value = 'feet'
$ws.Cells.Item(124,6).Value = 'They never age, they never wear out, they always look good. Always'
$ws.Cells.Item(124,7).Value = 'feet'
$ws.Cells.Item(124,10).Value = 22
$ws.Cells.Item(124,11).Value = 2500
$ws.Cells.Item(124,15).Value = 0.24
$ws.Cells.Item(124,16).Value = 0.24
$ws.Cells.Item(124,17).Value = 0.25
$ws.Cells.Item(124,18).Value = 0.24
$ws.Cells.Item(124,19).Value = 0.24
$ws.Cells.Item(124,21).Value = 1
$ws.Cells.Item(124,22).Value = 35
$ws.Cells.Item(124,23).Value = 78
$ws.Cells.Item(124,24).Value = 'armour'
$ws.Cells.Item(125,4).Value = 'Mythril sleeves'
$ws.Cells.Item(125,5).Value = 'sleeves'
$ws.Cells.Item(125,6).Value = 'Made from the rarest substance on earth, lighter then any feather and stronger then any dragon scale, mythril is an interesting subject to work with.'
$ws.Cells.Item(125,7).Value = 'sleeves'
$ws.Cells.Item(125,10).Value = 20
$ws.Cells.Item(125,11).Value = 3200
$ws.Cells.Item(125,15).Value = 0.23
$ws.Cells.Item(125,18).Value = 0.24
$ws.Cells.Item(125,21).Value = 1
$ws.Cells.Item(125,22).Value = 37
$ws.Cells.Item(125,23).Value = 80
$ws.Cells.Item(125,24).Value = 'armour'
$ws.Cells.Item(126,4).Value = 'Obsidian Helm'
$ws.Cells.Item(126,5).Value = 'helmet'
$ws.Cells.Item(126,6).Value = 'made from this glass like substance, it''s harder then steel. While it might look delicate I assure you it is not.'
$ws.Cells.Item(126,7).Value = 'helmet'
$ws.Cells.Item(126,11).Value = 3000
$ws.Cells.Item(126,16).Value = 0.25
$ws.Cells.Item(126,17).Value = 0.24
$ws.Cells.Item(127,1).Value = 'Archers Bane'
$ws.Cells.Item(127,3).Value = 1
$ws.Cells.Item(127,4).Value = 'Magma Helm'
$ws.Cells.Item(127,5).Value = 'helmet'
$ws.Cells.Item(127,6).Value = 'Crafted from enchanted magma, don''t worry it won''t burn....Or it...Shouldn''t'
$ws.Cells.Item(127,7).Value = 'helmet'
$ws.Cells.Item(127,10).Value = 20
$ws.Cells.Item(127,11).Value = 865
$ws.Cells.Item(127,15).Value = 0.17
$ws.Cells.Item(127,16).Value = 0.17
$ws.Cells.Item(127,19).Value = 0.17
$ws.Cells.Item(127,21).Value = 1
$ws.Cells.Item(127,22).Value = 31
$ws.Cells.Item(127,23).Value = 60
$ws.Cells.Item(127,24).Value = 'armour'
$ws.Cells.Item(128,4).Value = 'Angelic Steel Gloves'
$ws.Cells.Item(128,5).Value = 'gloves'
$ws.Cells.Item(128,6).Value = 'Angelic Steel Gloves come from the heavens above. The gates have always seemed closed and no one has ever been able to enter the heavenly plane. Maybe you can? Maybe they will open the doors for you?'
$ws.Cells.Item(128,7).Value = 'gloves'
$ws.Cells.Item(128,10).Value = 20
$ws.Cells.Item(128,11).Value = 1890
$ws.Cells.Item(128,15).Value = 0.2
$ws.Cells.Item(128,16).Value = 0.25
$ws.Cells.Item(128,17).Value = 0.2
$ws.Cells.Item(128,18).Value = 0.24
$ws.Cells.Item(128,19).Value = 0.24
$ws.Cells.Item(128,21).Value = 1
$ws.Cells.Item(128,22).Value = 36
$ws.Cells.Item(128,23).Value = 80
$ws.Cells.Item(128,24).Value = 'armour'
$ws.Cells.Item(129,4).Value = 'Crystal Ring'
$ws.Cells.Item(129,5).Value = 'ring'
$ws.Cells.Item(129,6).Value = 'Made completely from the rarest crystals on this plane, even the band is made from crystals. There is magic that radiates in this ring.'
$ws.Cells.Item(129,8).Value = 15
$ws.Cells.Item(129,11).Value = 5000
$ws.Cells.Item(129,15).Value = 0.24
$ws.Cells.Item(129,16).Value = 0.24
$ws.Cells.Item(129,17).Value = 0.24
$ws.Cells.Item(129,18).Value = 0.24
$ws.Cells.Item(129,19).Value = 0.24
$ws.Cells.Item(129,21).Value = 1
$ws.Cells.Item(129,22).Value = 25
$ws.Cells.Item(129,23).Value = 58
$ws.Cells.Item(129,24).Value = 'ring'
$ws.Cells.Item(130,4).Value = 'Hellhounds'
$ws.Cells.Item(130,5).Value = 'spell-damage'
$ws.Cells.Item(130,6).Value = 'Conjure the hounds of hell to devour your enemies.'
$ws.Cells.Item(130,8).Value = 185
$ws.Cells.Item(130,11).Value = 5600
$ws.Cells.Item(130,21).Value = 1
$ws.Cells.Item(130,22).Value = 16
$ws.Cells.Item(130,23).Value = 30
$ws.Cells.Item(130,24).Value = 'spell'
$ws.Cells.Item(131,4).Value = 'Kiss For The Reaper'
$ws.Cells.Item(131,5).Value = 'spell-healing'
$ws.Cells.Item(131,6).Value = 'The reaper want''s a kiss, a single kiss. Give the reaper a kiss and your wounds will fade away.'
$ws.Cells.Item(131,9).Value = 160
$ws.Cells.Item(131,11).Value = 6000
$ws.Cells.Item(131,21).Value = 1
$ws.Cells.Item(131,22).Value = 20
$ws.Cells.Item(131,23).Value = 40
$ws.Cells.Item(131,24).Value = 'spell'
$ws.Cells.Item(132,4).Value = 'Hateful Wish'
$ws.Cells.Item(132,5).Value = 'artifact'
$ws.Cells.Item(132,6).Value = 'Once it was made, the wish, hate was said to fill the land of the hearts o children. Children who rose up and slaughtered their parents in their beds. Alas fairy tales tend to get dark.'
$ws.Cells.Item(132,8).Value = 125
$ws.Cells.Item(132,9).Value = 15
$ws.Cells.Item(132,10).Value = 25
$ws.Cells.Item(132,11).Value = 3000
$ws.Cells.Item(132,18).Value = 0.25
$ws.Cells.Item(132,19).Value = 0.25
$ws.Cells.Item(132,21).Value = 1
$ws.Cells.Item(132,22).Value = 10
$ws.Cells.Item(132,23).Value = 30
$ws.Cells.Item(132,24).Value = 'artifact'

# Re-fit the description (F) and skill_level_required (K) columns now that
# they contain longer values than before.
$ws.Columns.Item(6).ColumnWidth = 250.25
$ws.Columns.Item(11).ColumnWidth = 6.2
